# [Kadastro App] Yeni kayit eklendi: 3003
# Appends a new record row (row 63) to both the master "Kayitlar" sheet and
# the filtered "Erdemli" sheet, mirroring the existing data layout where every
# column is stored as text (even numeric-looking values like the record id).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Kayitlar", "Erdemli")

$newRow = 63

$values = @("3003", "2025-09-11", "Erdemli", "1", "ÇAP", "AYHAN KARADAYI (K.Teknisyeni)")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $rowRange = $ws.Range("A" + $newRow + ":F" + $newRow)
    # Force text storage so values like "3003" / dates / "1" aren't
    # reinterpreted as numbers or dates, matching the rest of the sheet.
    $rowRange.NumberFormat = "@"

    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($newRow, $col).Value = $values[$col - 1]
    }
}
